$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct/recode a few values in the lookup table
$ws.Range("D7").Value = "Alameda"
$ws.Range("C7").Value = "Pensylvania"
$ws.Range("B5").Value = "USA"

# Update the selected cell to match the saved view state
$ws.Range("E7").Select()
